# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# for the cryptos worksheet, matching the commit's refreshed scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $rng = $ws.Range($CellRef)
    # Force text interpretation so numeric-looking strings (e.g. "0.530")
    # are not silently coerced into Double values and lose formatting,
    # then restore the default 'Normal' style so no stray number-format
    # style id gets attached to the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '67.290.51'
Set-TextCell 'D3' '2.546.63'
Set-TextCell 'E3' '  -2.69%  '
Set-TextCell 'D4' '0.999'
Set-TextCell 'E4' '  -0.05%  '
Set-TextCell 'D5' '591.48'
Set-TextCell 'E5' '  +0.29%  '
Set-TextCell 'D6' '174.79'
Set-TextCell 'E6' '  +5.56%  '
Set-TextCell 'E7' '  -0.02%  '
Set-TextCell 'D8' '0.530'
Set-TextCell 'E8' '  -0.34%  '
Set-TextCell 'D9' '2.545.14'
Set-TextCell 'E9' '  -2.74%  '
Set-TextCell 'D10' '0.140'
Set-TextCell 'E10' '  +1.24%  '
Set-TextCell 'E11' '  +1.17%  '
Set-TextCell 'E13' '  -0.72%  '
Set-TextCell 'E14' '  -0.99%  '
Set-TextCell 'D15' '3.012.38'
Set-TextCell 'E15' '  -2.58%  '
Set-TextCell 'E16' '  -0.79%  '
Set-TextCell 'D17' '67.110.35'
Set-TextCell 'E17' '  +0.15%  '
Set-TextCell 'D18' '2.536.11'
Set-TextCell 'E18' '  -3.16%  '
Set-TextCell 'D19' '8.07'
Set-TextCell 'E19' '  +3.26%  '
Set-TextCell 'E20' '  -3.19%  '
Set-TextCell 'D21' '355.85'
Set-TextCell 'E21' '  +0.27%  '
Set-TextCell 'E22' '  -1.24%  '
Set-TextCell 'E23' '  +0.90%  '
Set-TextCell 'E24' '  +4.12%  '
Set-TextCell 'E25' '  +0.03%  '
Set-TextCell 'D26' '70.03'
Set-TextCell 'E26' '  +1.34%  '
Set-TextCell 'D27' '10.08'
Set-TextCell 'E27' '  -4.18%  '
Set-TextCell 'D28' '2.685.13'
Set-TextCell 'E28' '  -2.37%  '
Set-TextCell 'E29' '  +0.05%  '
Set-TextCell 'D30' '0.0₃0997'
Set-TextCell 'E30' '  +0.01%  '
Set-TextCell 'D31' '536.21'
Set-TextCell 'E31' '  -1.35%  '
Set-TextCell 'D32' '8.27'
Set-TextCell 'E32' '  +5.09%  '
Set-TextCell 'E33' '  +0.82%  '
Set-TextCell 'E34' '  -0.54%  '
Set-TextCell 'E35' '  -1.52%  '
Set-TextCell 'D36' '0.999'
Set-TextCell 'E36' '  -0.03%  '
Set-TextCell 'E37' '  -0.40%  '
Set-TextCell 'D38' '157.59'
Set-TextCell 'E38' '  +0.30%  '
Set-TextCell 'E39' '  -0.68%  '
Set-TextCell 'E40' '  +1.10%  '
Set-TextCell 'E41' '  -1.98%  '
Set-TextCell 'E42' '  +0.27%  '
Set-TextCell 'D43' '5.21'
Set-TextCell 'E43' '  +1.38%  '
Set-TextCell 'E44' '  +6.71%  '
Set-TextCell 'E45' '  -0.02%  '
Set-TextCell 'D46' '39.84'
Set-TextCell 'E46' '  -0.99%  '
Set-TextCell 'D47' '151.48'
Set-TextCell 'E47' '  +0.07%  '
Set-TextCell 'E49' '  -6.18%  '
Set-TextCell 'E50' '  -1.29%  '
Set-TextCell 'E51' '  +0.95%  '
